$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the timing values with new, more precise measurements.
# Filled column by column (B then C then D) to match the shared-string order
# produced by the original authoring session.

# Column B (1000 operations)
$ws.Range("B2").Value = "0.04363 ms"
$ws.Range("B3").Value = "0.00015 ms"
$ws.Range("B4").Value = "0.00016 ms"
$ws.Range("B5").Value = "0.01130 ms"
$ws.Range("B6").Value = "0.48770 ms"
$ws.Range("B7").Value = "0.02550 ms"

# Column C (10000 operations)
$ws.Range("C2").Value = "0.44977 ms"
$ws.Range("C3").Value = "0.00467 ms"
$ws.Range("C4").Value = "0.00055 ms"
$ws.Range("C5").Value = "0.02520 ms"
$ws.Range("C6").Value = "0.36920 ms"
$ws.Range("C7").Value = "0.03630 ms"

# Column D (100000 operations)
$ws.Range("D2").Value = "6.06718 ms"
$ws.Range("D3").Value = "0.00046 ms"
$ws.Range("D4").Value = "0.00032 ms"
$ws.Range("D5").Value = "0.24200 ms"
$ws.Range("D6").Value = "36.78780 ms"
$ws.Range("D7").Value = "0.53920 ms"

# Move the active selection to D7, matching the saved view state.
$ws.Range("D7").Select()
